$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Summary sheet - update aggregate metrics after trade #21 closed
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B5").Value = -0.2     # Total P&L %
$summary.Range("B6").Value = 21       # Total Trades
$summary.Range("B9").Value = 47.62    # Win Rate %

# ---------------------------------------------------------------------------
# 2) Strategy Status sheet - MarketMaking row (row 5) reflects the new trade
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D5").Value = 21        # Trades
$status.Range("G5").Value = 47.62     # Win Rate %

# ---------------------------------------------------------------------------
# Helper: append the new trade-log row (row 22) to a trade-log sheet.
# Columns B (date) and C (time) hold text that LOOKS like a date/time
# ("2026-02-17"); a plain .Value assignment of the date-like string gets
# silently reinterpreted by Excel as a date serial number. Pre-formatting
# the cell as Text ("@") keeps the literal string, and resetting the style
# back to "Normal" afterwards (without touching .Value again) keeps the
# cell on the default style, matching a cell that was never reformatted.
# ---------------------------------------------------------------------------
function Add-TradeRow22 {
    param($ws)

    $ws.Cells.Item(22, 1).Value = 21

    $ws.Cells.Item(22, 2).NumberFormat = "@"
    $ws.Cells.Item(22, 2).Value = "2026-02-17"
    $ws.Cells.Item(22, 2).Style = "Normal"

    $ws.Cells.Item(22, 3).Value = "20:04:34"

    $ws.Cells.Item(22, 4).Value = "MarketMaking"
    $ws.Cells.Item(22, 5).Value = "DOWN"
    $ws.Cells.Item(22, 6).Value = 0.02
    $ws.Cells.Item(22, 7).Value = 0.02
    $ws.Cells.Item(22, 8).Value = "CLOSED"
    $ws.Cells.Item(22, 9).Value = 0
    $ws.Cells.Item(22, 10).Value = 0
    $ws.Cells.Item(22, 11).Value = 99.8
    $ws.Cells.Item(22, 12).Value = 0
    $ws.Cells.Item(22, 13).Value = 0
    $ws.Cells.Item(22, 14).Value = 0.6
    $ws.Cells.Item(22, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(22, 16).Value = "early_exit"
    $ws.Cells.Item(22, 17).Value = 0.13
}

# ---------------------------------------------------------------------------
# 3) All Trades sheet - append trade #21 as row 22
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow22 $allTrades

# ---------------------------------------------------------------------------
# 4) MarketMaking sheet - same trade log, append trade #21 as row 22
# ---------------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow22 $marketMaking
